$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 35-36 (a new weekly price report for Perejil,
# date 45072), shifting the existing rows 35-65 down to 37-67.
$ws.Rows("35:36").Insert()

# Row 35: Perejil, Primera, date 45072
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 45072
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112044
$ws.Range("G35").Value = "Perejil"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 1200
$ws.Range("L35").Value = 1200
$ws.Range("M35").Value = 1200
$ws.Range("N35").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 1200
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"

# Row 36: Perejil, Segunda, date 45072
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 45072
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112044
$ws.Range("G36").Value = "Perejil"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = 1000
$ws.Range("N36").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 1000
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"
